# 2050_TM152_DBP_NoProject_04 is complete and now current 2050 noproject
#
# Row 34 (2050_TM152_DBP_NoProject_01) was "current" -> no longer current.
# Row 35 (2050_TM152_DBP_NoProject_03) was "running" -> no longer running.
# A new row is inserted as row 36 for 2050_TM152_DBP_NoProject_04, copying
# the urbansim_path/urbansim_runid from the (formerly running) row and
# marking it "current" - the new run has completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "current"/"running" status markers on the two prior NoProject rows.
$ws.Range("H34").Value = $null
$ws.Range("H35").Value = $null

# Insert a new row at 36, pushing the existing rows 36-42 down to 37-43.
$ws.Rows("36:36").Insert(-4121)

# Copy formatting from the row right below (same style block) onto the new row.
$ws.Range("A37:H37").Copy()
$ws.Range("A36:H36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A36").Value = "RTP2021"
$ws.Range("B36").Value = 2050
$ws.Range("C36").Value = "2050_TM152_DBP_NoProject_04"
$ws.Range("D36").Value = "DraftBlueprint"
$ws.Range("E36").Value = "No Project"
$ws.Range("F36").Value = """Blueprint Plus Crossing (s23)\v1.5.2"""
$ws.Range("G36").Value = "run72"
$ws.Range("H36").Value = "current"

# Match the author's final cursor position recorded in the saved file.
[void]$ws.Range("C40").Select()
